$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume(1h) (E) columns for each coin row
# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.852.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.35%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.657.75"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.18%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.54%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.05%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3631"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.82%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.17"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.85%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3251"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.22%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.132"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.95%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07053"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.40%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.05%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.032"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.97%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.50"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.39%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.661.44"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.48%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.579"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.02%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001044"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.66%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06590"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.03%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.03%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "78.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.69%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.894"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.18%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.74"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.19%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.77%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.847.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.18%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.424"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.35%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.451"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.90%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "147.20"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.75%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.57"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.04%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.843.03"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.70%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.080"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.63%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.757"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -10.45%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08462"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.42%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.652"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.41%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.25"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.03%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.285"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.89%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.162"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.36%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02257"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.16%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06048"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.94%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.333"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.71%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2067"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.85%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9990"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.07%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5927"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.98%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.770"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.44%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.72"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.20%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5614"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.54%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.93"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.05%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.943"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.32%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06967"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.72%  "

# Row 51
$ws.Range("E51").Value = "  -0.24%  "

# Rows 30 and 31 swapped: BitcoinCash <-> ImmutableX (with updated data)
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.209"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.55%  "

$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "124.83"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.00%  "
